$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7496711949059137
$ws.Range("C2").Value = -0.1800508050940863
$ws.Range("D2").Value = 0.05495719490591373
$ws.Range("E2").Value = -0.2521788050940862
$ws.Range("F2").Value = 0.2584971949059137
$ws.Range("G2").Value = -0.02886980509408627
$ws.Range("H2").Value = -0.0008708050940862688
$ws.Range("I2").Value = 0.09177719490591374
$ws.Range("J2").Value = 0.3047531949059137
$ws.Range("K2").Value = -0.1959068050940863
$ws.Range("B3").Value = -0.1976049264540507
$ws.Range("C3").Value = 0.03740307354594928
$ws.Range("D3").Value = -0.2697329264540507
$ws.Range("E3").Value = 0.2409430735459493
$ws.Range("F3").Value = -0.04642392645405072
$ws.Range("G3").Value = -0.01842492645405072
$ws.Range("H3").Value = 0.07422307354594927
$ws.Range("I3").Value = 0.2871990735459493
$ws.Range("J3").Value = -0.2134609264540507
$ws.Range("K3").Value = 0.07629607354594928
$ws.Range("B4").Value = 0.1255103924969555
$ws.Range("C4").Value = -0.1816256075030445
$ws.Range("D4").Value = 0.3290503924969555
$ws.Range("E4").Value = 0.0416833924969555
$ws.Range("F4").Value = 0.06968239249695551
$ws.Range("G4").Value = 0.1623303924969555
$ws.Range("H4").Value = 0.3753063924969555
$ws.Range("I4").Value = -0.1253536075030445
$ws.Range("J4").Value = 0.1644033924969555
$ws.Range("K4").Value = -0.0934886075030445
$ws.Range("B5").Value = -0.003050974015260888
$ws.Range("C5").Value = 0.5076250259847391
$ws.Range("D5").Value = 0.2202580259847391
$ws.Range("E5").Value = 0.2482570259847391
$ws.Range("F5").Value = 0.3409050259847391
$ws.Range("G5").Value = 0.5538810259847391
$ws.Range("H5").Value = 0.0532210259847391
$ws.Range("I5").Value = 0.3429780259847391
$ws.Range("J5").Value = 0.08508602598473911
$ws.Range("K5").Value = 0.4674670259847391
$ws.Range("B6").Value = 1.419652253737389
$ws.Range("C6").Value = 1.132285253737389
$ws.Range("D6").Value = 1.160284253737389
$ws.Range("E6").Value = 1.252932253737389
$ws.Range("F6").Value = 1.465908253737389
$ws.Range("G6").Value = 0.9652482537373891
$ws.Range("H6").Value = 1.255005253737389
$ws.Range("I6").Value = 0.9971132537373891
$ws.Range("J6").Value = 1.379494253737389
$ws.Range("K6").Value = 0.9682171597373891
$ws.Range("B7").Value = 0.1883110177716323
$ws.Range("C7").Value = 0.2163100177716323
$ws.Range("D7").Value = 0.3089580177716323
$ws.Range("E7").Value = 0.5219340177716323
$ws.Range("F7").Value = 0.02127401777163229
$ws.Range("G7").Value = 0.3110310177716323
$ws.Range("H7").Value = 0.0531390177716323
$ws.Range("I7").Value = 0.4355200177716323
$ws.Range("J7").Value = 0.02424292377163231
$ws.Range("B8").Value = 0.2202779152847414
$ws.Range("C8").Value = 0.3129259152847414
$ws.Range("D8").Value = 0.5259019152847414
$ws.Range("E8").Value = 0.0252419152847414
$ws.Range("F8").Value = 0.3149989152847414
$ws.Range("G8").Value = 0.05710691528474141
$ws.Range("H8").Value = 0.4394879152847414
$ws.Range("I8").Value = 0.02821082128474142
$ws.Range("B9").Value = 0.4485660054549828
$ws.Range("C9").Value = 0.6615420054549828
$ws.Range("D9").Value = 0.1608820054549828
$ws.Range("E9").Value = 0.4506390054549828
$ws.Range("F9").Value = 0.1927470054549828
$ws.Range("G9").Value = 0.5751280054549828
$ws.Range("H9").Value = 0.1638509114549828
$ws.Range("B10").Value = 0.420735823599318
$ws.Range("C10").Value = -0.07992417640068192
$ws.Range("D10").Value = 0.2098328235993181
$ws.Range("E10").Value = -0.04805917640068191
$ws.Range("F10").Value = 0.3343218235993181
$ws.Range("G10").Value = -0.0769552704006819
$ws.Range("B11").Value = -0.1252553916527783
$ws.Range("C11").Value = 0.1645016083472217
$ws.Range("D11").Value = -0.09339039165277828
$ws.Range("E11").Value = 0.2889906083472217
$ws.Range("F11").Value = -0.1222864856527783
$ws.Range("B12").Value = 0.09764018641116785
$ws.Range("C12").Value = -0.1602518135888321
$ws.Range("D12").Value = 0.2221291864111679
$ws.Range("E12").Value = -0.1891479075888321
$ws.Range("B13").Value = -0.1898380159455487
$ws.Range("C13").Value = 0.1925429840544512
$ws.Range("D13").Value = -0.2187341099455487
$ws.Range("B14").Value = 0.1743923273248104
$ws.Range("C14").Value = -0.2368847666751896
$ws.Range("B15").Value = -0.254916590923889
$ws.Range("B16").Value = -0.01954473626955232
